$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Selene " / "Anahí" / " POSADAS" were three separate runs; collapse
#    them into a single run with the same text. A Find/Replace across the
#    run boundary is the simplest way to get Word to re-merge them into
#    one run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Selene Anahí POSADAS", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Selene Anahí POSADAS", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Append the "ACTUALIZACION 18/11/2021" section right after the
#    "Tamara ROMERO" paragraph (the last entry of "Integrantes del
#    equipo"), and move the hidden "_GoBack" bookmark from the end of
#    "Tamara ROMERO" to the end of the newly-added content.
# ---------------------------------------------------------------------

# Locate the "Tamara ROMERO" paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Tamara ROMERO") {
        $target = $p
    }
}

# Drop the existing hidden bookmark; the OOXML fragment below re-creates
# it in its new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$countBefore = $d.Paragraphs.Count

# OOXML (WordprocessingML wrapped in the pkg:package form Range.InsertXML
# expects) describing the six new paragraphs to insert. A trailing bare
# <w:p/> is appended on purpose: Word always folds the *last* paragraph
# mark of inserted XML into whatever paragraph follows the insertion
# point, so without this spacer the real final paragraph ("El resto...")
# would inherit the formatting of -- and merge into -- the document's
# pre-existing trailing empty paragraph. The spacer absorbs that merge
# instead, and is deleted again afterwards.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>ACTUALIZACION 18/11/2021</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Profe, el motivo de esta actualización es que debíamos dejar el sitio funcionando para la materia de practica </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>profesionalizante</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>, dado el caso hemos realizado algunas modificaciones a lo anteriormente mencionado.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Para precisar, en el punto 3, ahora al momento de cargar los datos, el sistema nos informa que los datos se han cargado de forma correcta</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> y nos redirigirá a la página de login.html.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Al poner nuestro usuario y contraseña la página nos lleva a la anteriormente programada página “en construcción” nos saludará por nuestro nombre y apellido, pero ahora extrayendo los datos directo de la base de datos, luego de 8 segundos nos re direccionará a el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>index</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> como se encontraba realizado con anterioridad.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>El resto de los puntos se siguen manteniendo de la misma forma.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint = $d.Range($target.Range.End, $target.Range.End)
$insertPoint.InsertXML($xml) | Out-Null

# Remove the spacer paragraph that absorbed the merge-with-next-paragraph
# behaviour described above. It sits right after the last of the six new
# paragraphs, i.e. immediately before the document's original trailing
# (empty, underlined) paragraph.
$spacerIndex = $countBefore + 6
$spacer = $d.Paragraphs.Item($spacerIndex)
$spacer.Range.Delete()

Write-Output "done"
